# Revision de no conformidades
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the real closing date for the first non-conformidad (row 4, column E)
$ws.Range("E4").Value = Get-Date -Year 2015 -Month 12 -Day 22 -Hour 0 -Minute 0 -Second 0

# Mark the status of that entry as "Cerrada" (was "En proceso")
$ws.Range("F4").Value = "Cerrada"

# Move the active selection to E5 (as reflected in the sheet view)
$ws.Range("E5").Select()
